$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 18
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03).
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45202
}
